$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 1.36
$ws.Range("G2").Value = 1.42
$ws.Range("H2").Value = 6.2
$ws.Range("I2").Value = 10.5
$ws.Range("J2").Value = 4.2
$ws.Range("K2").Value = 8
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 3.45
$ws.Range("O2").Value = 1.15
$ws.Range("P2").Value = 2.52
$ws.Range("Q2").Value = 1.43
$ws.Range("R2").Value = 1.7
$ws.Range("S2").Value = 2.06
$ws.Range("T2").Value = 1.68
$ws.Range("U2").Value = 2.12
$ws.Range("V2").Value = 1.1
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 40
$ws.Range("Y2").Value = 44
$ws.Range("Z2").Value = 95
$ws.Range("AB2").Value = 15.5
$ws.Range("AC2").Value = 16.5
$ws.Range("AD2").Value = 38
$ws.Range("AF2").Value = 13.5
$ws.Range("AG2").Value = 13
$ws.Range("AH2").Value = 27
$ws.Range("AI2").Value = 95
$ws.Range("AJ2").Value = 16
$ws.Range("AK2").Value = 17
$ws.Range("AL2").Value = 34
$ws.Range("AN2").Value = 5.4
$ws.Range("AO2").Value = 110
$ws.Range("F3").Value = 1.32
$ws.Range("G3").Value = 1.64
$ws.Range("J3").Value = 4.3
$ws.Range("K3").Value = 12
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 2.18
$ws.Range("O3").Value = 1.18
$ws.Range("P3").Value = 2.18
$ws.Range("Q3").Value = 1.43
$ws.Range("R3").Value = 1.47
$ws.Range("S3").Value = 2.08
$ws.Range("W3").Value = 2.54
$ws.Range("W4").Value = 1.14
$ws.Range("F5").Value = 3.3
$ws.Range("G5").Value = 5.5
$ws.Range("H5").Value = 1.01
$ws.Range("I5").Value = 2.16
$ws.Range("J5").Value = 3.75
$ws.Range("K5").Value = 950
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 2.2
$ws.Range("P5").Value = 2.2
$ws.Range("Q5").Value = 1.48
$ws.Range("R5").Value = 1.5
$ws.Range("V5").Value = 1.86
$ws.Range("W5").Value = 1.22
$ws.Range("G6").Value = 6.2
$ws.Range("K6").Value = 950
$ws.Range("W6").Value = 1.2
$ws.Range("F7").Value = 1.15
$ws.Range("G7").Value = 1.23
$ws.Range("H7").Value = 8
$ws.Range("I7").Value = 38
$ws.Range("K7").Value = 9.6
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 4.2
$ws.Range("O7").Value = 1.2
$ws.Range("P7").Value = 2.16
$ws.Range("Q7").Value = 1.58
$ws.Range("R7").Value = 1.54
$ws.Range("S7").Value = 2.26
$ws.Range("T7").Value = 2.68
$ws.Range("U7").Value = 1.47
$ws.Range("V7").Value = 1.02
$ws.Range("W7").Value = 5.1
$ws.Range("X7").Value = 32
$ws.Range("Y7").Value = 85
$ws.Range("Z7").Value = 390
$ws.Range("AB7").Value = 10.5
$ws.Range("AC7").Value = 25
$ws.Range("AD7").Value = 130
$ws.Range("AF7").Value = 8.2
$ws.Range("AG7").Value = 16.5
$ws.Range("AH7").Value = 75
$ws.Range("AI7").Value = 590
$ws.Range("AJ7").Value = 9.2
$ws.Range("AK7").Value = 20
$ws.Range("AL7").Value = 85
$ws.Range("AM7").Value = 550
$ws.Range("AN7").Value = 4.6
$ws.Range("G9").Value = 140
$ws.Range("H9").Value = 1.4
$ws.Range("I9").Value = 1.77
$ws.Range("K9").Value = 11
$ws.Range("G10").Value = 2.64
$ws.Range("H10").Value = 2.86
$ws.Range("I10").Value = 3.3
$ws.Range("J10").Value = 3.4
$ws.Range("R10").Value = 1.39
$ws.Range("T10").Value = 1.68
$ws.Range("U10").Value = 2.2
$ws.Range("V10").Value = 1.44
$ws.Range("W10").Value = 1.6
$ws.Range("Z10").Value = 26
$ws.Range("AE10").Value = 40
$ws.Range("X12").Value = 970
$ws.Range("Y12").Value = 970
$ws.Range("AB12").Value = 970
$ws.Range("AC12").Value = 970
$ws.Range("AF12").Value = 970
$ws.Range("AG12").Value = 970
$ws.Range("G13").Value = 1.6
$ws.Range("J14").Value = 1.03
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2.38
$ws.Range("H15").Value = 3.1
$ws.Range("I15").Value = 4.1
$ws.Range("J15").Value = 3.2
$ws.Range("K15").Value = 4
$ws.Range("N15").Value = 3.1
$ws.Range("O15").Value = 1.29
$ws.Range("Q15").Value = 1.78
$ws.Range("T15").Value = 1.72
$ws.Range("U15").Value = 2.06
$ws.Range("V15").Value = 1.32
$ws.Range("W15").Value = 1.72
$ws.Range("X15").Value = 18
$ws.Range("Y15").Value = 17
$ws.Range("Z15").Value = 32
$ws.Range("AA15").Value = 80
$ws.Range("AB15").Value = 12.5
$ws.Range("AD15").Value = 18
$ws.Range("AE15").Value = 55
$ws.Range("AF15").Value = 17.5
$ws.Range("AG15").Value = 13.5
$ws.Range("AH15").Value = 21
$ws.Range("AI15").Value = 60
$ws.Range("AJ15").Value = 32
$ws.Range("AK15").Value = 29
$ws.Range("AL15").Value = 46
$ws.Range("AM15").Value = 110
$ws.Range("AN15").Value = 21
$ws.Range("AO15").Value = 50
$ws.Range("F16").Value = 6.6
$ws.Range("G16").Value = 9.2
$ws.Range("H16").Value = 1.42
$ws.Range("I16").Value = 1.5
$ws.Range("J16").Value = 4.7
$ws.Range("K16").Value = 5.6
$ws.Range("P16").Value = 2.28
$ws.Range("Q16").Value = 1.63
$ws.Range("S16").Value = 2.6
$ws.Range("T16").Value = 1.87
$ws.Range("U16").Value = 1.94
$ws.Range("V16").Value = 3
$ws.Range("W16").Value = 1.12
$ws.Range("AA16").Value = 15.5
$ws.Range("AE16").Value = 18
$ws.Range("AF16").Value = 90
$ws.Range("AJ16").Value = 310
$ws.Range("AK16").Value = 150
$ws.Range("AM16").Value = 150
$ws.Range("AN16").Value = 170
$ws.Range("AO16").Value = 7.2
$ws.Range("G17").Value = 1.9
$ws.Range("L17").Value = 1.44
$ws.Range("R17").Value = 1.28
$ws.Range("W17").Value = 2.1
$ws.Range("Z17").Value = 42
$ws.Range("AB17").Value = 7.8
$ws.Range("AF17").Value = 11
$ws.Range("AG17").Value = 11
$ws.Range("AN17").Value = 16
$ws.Range("L18").Value = 1.32
$ws.Range("AI18").Value = 42
$ws.Range("F19").Value = 1.99
$ws.Range("I19").Value = 4.7
$ws.Range("J19").Value = 3.15
$ws.Range("K19").Value = 3.8
$ws.Range("L19").Value = 1.43
$ws.Range("M19").Value = 1.08
$ws.Range("P19").Value = 1.77
$ws.Range("Q19").Value = 1.97
$ws.Range("R19").Value = 1.29
$ws.Range("S19").Value = 3.7
$ws.Range("T19").Value = 1.83
$ws.Range("U19").Value = 1.96
$ws.Range("V19").Value = 1.28
$ws.Range("W19").Value = 1.78
$ws.Range("X19").Value = 970
$ws.Range("Y19").Value = 970
$ws.Range("Z19").Value = 34
$ws.Range("AA19").Value = 110
$ws.Range("AB19").Value = 9.8
$ws.Range("AC19").Value = 9
$ws.Range("AD19").Value = 970
$ws.Range("AE19").Value = 65
$ws.Range("AF19").Value = 970
$ws.Range("AG19").Value = 970
$ws.Range("AH19").Value = 23
$ws.Range("AI19").Value = 75
$ws.Range("AJ19").Value = 32
$ws.Range("AK19").Value = 28
$ws.Range("AL19").Value = 48
$ws.Range("AM19").Value = 140
$ws.Range("AN19").Value = 22
$ws.Range("AO19").Value = 75
$ws.Range("J21").Value = 3.45
$ws.Range("L21").Value = 1.32
$ws.Range("W22").Value = 1.75
$ws.Range("Y23").Value = 17
$ws.Range("AB23").Value = 6
$ws.Range("AD23").Value = 32
$ws.Range("G24").Value = 3.5
$ws.Range("H24").Value = 2.2
$ws.Range("I24").Value = 2.22
$ws.Range("R24").Value = 1.51
$ws.Range("S24").Value = 2.88
$ws.Range("V24").Value = 1.81
$ws.Range("Z24").Value = 15.5
$ws.Range("AF24").Value = 26
$ws.Range("AN24").Value = 27
$ws.Range("G25").Value = 3.05
$ws.Range("K25").Value = 3.95
$ws.Range("P25").Value = 2.42
$ws.Range("Q25").Value = 1.68
$ws.Range("R25").Value = 1.57
$ws.Range("W25").Value = 1.48
$ws.Range("H26").Value = 26
$ws.Range("I26").Value = 28
$ws.Range("Q26").Value = 1.37
$ws.Range("S26").Value = 1.91
$ws.Range("T26").Value = 2.38
$ws.Range("AE26").Value = 510
$ws.Range("AI26").Value = 300
$ws.Range("AL26").Value = 46
$ws.Range("AO26").Value = 610
$ws.Range("F27").Value = 1.71
$ws.Range("G27").Value = 1.72
$ws.Range("H27").Value = 5.2
$ws.Range("I27").Value = 5.3
$ws.Range("L27").Value = 1.29
$ws.Range("P27").Value = 2.7
$ws.Range("S27").Value = 2.4
$ws.Range("U27").Value = 2.48
$ws.Range("V27").Value = 1.23
$ws.Range("W27").Value = 2.38
$ws.Range("AD27").Value = 19
$ws.Range("AF27").Value = 12.5
$ws.Range("F28").Value = 3.05
$ws.Range("G28").Value = 3.1
$ws.Range("H28").Value = 2.5
$ws.Range("I28").Value = 2.52
$ws.Range("J28").Value = 3.65
$ws.Range("Q28").Value = 1.75
$ws.Range("T28").Value = 1.62
$ws.Range("V28").Value = 1.66
$ws.Range("Y28").Value = 13
$ws.Range("Z28").Value = 17
$ws.Range("AE28").Value = 23
$ws.Range("AF28").Value = 22
$ws.Range("AJ28").Value = 48
$ws.Range("AL28").Value = 36
$ws.Range("F29").Value = 1.25
$ws.Range("G29").Value = 1.26
$ws.Range("H29").Value = 13
$ws.Range("I29").Value = 13.5
$ws.Range("J29").Value = 7.8
$ws.Range("K29").Value = 8
$ws.Range("M29").Value = 1.01
$ws.Range("N29").Value = 9.6
$ws.Range("R29").Value = 2.16
$ws.Range("T29").Value = 1.79
$ws.Range("U29").Value = 2.22
$ws.Range("V29").Value = 1.08
$ws.Range("W29").Value = 4.9
$ws.Range("Z29").Value = 140
$ws.Range("AD29").Value = 46
$ws.Range("AE29").Value = 150
$ws.Range("AI29").Value = 110
$ws.Range("AJ29").Value = 11.5
$ws.Range("AM29").Value = 100
$ws.Range("AN29").Value = 3
$ws.Range("F30").Value = 5.4
$ws.Range("H30").Value = 1.63
$ws.Range("I30").Value = 1.64
$ws.Range("J30").Value = 4.9
$ws.Range("K30").Value = 5
$ws.Range("Q30").Value = 1.37
$ws.Range("S30").Value = 1.92
$ws.Range("V30").Value = 2.56
$ws.Range("X30").Value = 40
$ws.Range("AB30").Value = 40
$ws.Range("AF30").Value = 60
$ws.Range("AG30").Value = 22
$ws.Range("AJ30").Value = 130
$ws.Range("AN30").Value = 29
$ws.Range("AO30").Value = 4.8
$ws.Range("H31").Value = 2.46
$ws.Range("I31").Value = 2.5
$ws.Range("K31").Value = 4.1
$ws.Range("P31").Value = 2.64
$ws.Range("T31").Value = 1.53
$ws.Range("V31").Value = 1.67
$ws.Range("Y31").Value = 16
$ws.Range("AA31").Value = 36
$ws.Range("AB31").Value = 18
$ws.Range("AH31").Value = 13.5
$ws.Range("AK31").Value = 26
$ws.Range("AN31").Value = 16.5
$ws.Range("F32").Value = 3.55
$ws.Range("G32").Value = 3.6
$ws.Range("H32").Value = 2.26
$ws.Range("I32").Value = 2.28
$ws.Range("K32").Value = 3.6
$ws.Range("O32").Value = 1.25
$ws.Range("P32").Value = 2.24
$ws.Range("R32").Value = 1.49
$ws.Range("V32").Value = 1.78
$ws.Range("W32").Value = 1.38
$ws.Range("AA32").Value = 28
$ws.Range("AE32").Value = 21
$ws.Range("AL32").Value = 42
$ws.Range("AO32").Value = 14
$ws.Range("J33").Value = 4.4
$ws.Range("K33").Value = 4.5
$ws.Range("M33").Value = 1.03
$ws.Range("R33").Value = 1.74
$ws.Range("U33").Value = 2.66
$ws.Range("X33").Value = 27
$ws.Range("AH33").Value = 15.5
$ws.Range("F34").Value = 3.45
$ws.Range("G34").Value = 3.5
$ws.Range("H34").Value = 2.24
$ws.Range("I34").Value = 2.26
$ws.Range("N34").Value = 4.4
$ws.Range("P34").Value = 2.16
$ws.Range("Q34").Value = 1.84
$ws.Range("S34").Value = 3.15
$ws.Range("T34").Value = 1.71
$ws.Range("V34").Value = 1.79
$ws.Range("W34").Value = 1.4
$ws.Range("X34").Value = 16.5
$ws.Range("AB34").Value = 15
$ws.Range("AC34").Value = 8.2
$ws.Range("AE34").Value = 21
$ws.Range("F35").Value = 2.62
$ws.Range("G35").Value = 2.66
$ws.Range("P35").Value = 1.96
$ws.Range("S35").Value = 3.5
$ws.Range("V35").Value = 1.49
$ws.Range("W35").Value = 1.6
$ws.Range("Z35").Value = 19.5
$ws.Range("AF35").Value = 16
$ws.Range("AG35").Value = 13
$ws.Range("AI35").Value = 44
$ws.Range("AJ35").Value = 38
$ws.Range("AN35").Value = 22
$ws.Range("F36").Value = 1.55
$ws.Range("J36").Value = 5
$ws.Range("K36").Value = 5.1
$ws.Range("L36").Value = 1.27
$ws.Range("N36").Value = 6.4
$ws.Range("P36").Value = 2.8
$ws.Range("Q36").Value = 1.53
$ws.Range("R36").Value = 1.73
$ws.Range("S36").Value = 2.32
$ws.Range("T36").Value = 1.66
$ws.Range("W36").Value = 2.78
$ws.Range("Y36").Value = 32
$ws.Range("AC36").Value = 11.5
$ws.Range("AF36").Value = 11.5
$ws.Range("AJ36").Value = 15
$ws.Range("AK36").Value = 13
$ws.Range("AN36").Value = 5.7
$ws.Range("F37").Value = 2.58
$ws.Range("G37").Value = 2.6
$ws.Range("H37").Value = 3
$ws.Range("I37").Value = 3.05
$ws.Range("J37").Value = 3.55
$ws.Range("K37").Value = 3.6
$ws.Range("L37").Value = 1.4
$ws.Range("O37").Value = 1.28
$ws.Range("P37").Value = 2.14
$ws.Range("Q37").Value = 1.85
$ws.Range("R37").Value = 1.45
$ws.Range("S37").Value = 3.15
$ws.Range("U37").Value = 2.38
$ws.Range("V37").Value = 1.49
$ws.Range("W37").Value = 1.62
$ws.Range("X37").Value = 16
$ws.Range("Y37").Value = 13.5
$ws.Range("AF37").Value = 17
$ws.Range("AG37").Value = 11.5
$ws.Range("AJ37").Value = 36
$ws.Range("AK37").Value = 25
$ws.Range("AN37").Value = 18.5
$ws.Range("F38").Value = 1.34
$ws.Range("G38").Value = 1.35
$ws.Range("H38").Value = 9.4
$ws.Range("I38").Value = 9.8
$ws.Range("J38").Value = 6.6
$ws.Range("K38").Value = 6.8
$ws.Range("L38").Value = 1.2
$ws.Range("M38").Value = 1.01
$ws.Range("N38").Value = 9.4
$ws.Range("O38").Value = 1.11
$ws.Range("P38").Value = 3.85
$ws.Range("Q38").Value = 1.33
$ws.Range("R38").Value = 2.18
$ws.Range("S38").Value = 1.82
$ws.Range("T38").Value = 1.59
$ws.Range("U38").Value = 2.56
$ws.Range("V38").Value = 1.11
$ws.Range("W38").Value = 3.85
$ws.Range("X38").Value = 60
$ws.Range("Z38").Value = 110
$ws.Range("AA38").Value = 290
$ws.Range("AB38").Value = 16.5
$ws.Range("AC38").Value = 16.5
$ws.Range("AD38").Value = 36
$ws.Range("AE38").Value = 100
$ws.Range("AF38").Value = 12.5
$ws.Range("AG38").Value = 10.5
$ws.Range("AH38").Value = 21
$ws.Range("AI38").Value = 75
$ws.Range("AJ38").Value = 14
$ws.Range("AL38").Value = 19.5
$ws.Range("AM38").Value = 75
$ws.Range("AN38").Value = 3.45
$ws.Range("AO38").Value = 65
$ws.Range("I39").Value = 70
$ws.Range("J39").Value = 22
$ws.Range("T39").Value = 2.74
$ws.Range("U39").Value = 1.52
$ws.Range("W39").Value = 16
$ws.Range("AC39").Value = 990
$ws.Range("AG39").Value = 25
$ws.Range("AK39").Value = 18.5
$ws.Range("AL39").Value = 90
$ws.Range("AM39").Value = 530
$ws.Range("H40").Value = 1.87
$ws.Range("N40").Value = 5.2
$ws.Range("R40").Value = 1.57
$ws.Range("T40").Value = 1.64
$ws.Range("X40").Value = 21
$ws.Range("AC40").Value = 9.4
$ws.Range("AI40").Value = 26
$ws.Range("AJ40").Value = 90
$ws.Range("G41").Value = 2.58
$ws.Range("H41").Value = 3.85
$ws.Range("K41").Value = 2.96
$ws.Range("N41").Value = 2.34
$ws.Range("T41").Value = 2.36
$ws.Range("W41").Value = 1.64
$ws.Range("Y41").Value = 9.6
$ws.Range("Z41").Value = 36
$ws.Range("AB41").Value = 7.2
$ws.Range("AC41").Value = 6.8
$ws.Range("AD41").Value = 27
$ws.Range("AF41").Value = 13.5
$ws.Range("AG41").Value = 17
$ws.Range("AJ41").Value = 70
$ws.Range("AK41").Value = 100
$ws.Range("AN41").Value = 110
$ws.Range("F42").Value = 1.14
$ws.Range("G42").Value = 1.15
$ws.Range("I42").Value = 25
$ws.Range("J42").Value = 11
$ws.Range("K42").Value = 11.5
$ws.Range("N42").Value = 13
$ws.Range("R42").Value = 2.56
$ws.Range("W42").Value = 7.6
$ws.Range("Z42").Value = 290
$ws.Range("AB42").Value = 22
$ws.Range("AC42").Value = 27
$ws.Range("AD42").Value = 80
$ws.Range("AE42").Value = 320
$ws.Range("AF42").Value = 12.5
$ws.Range("AH42").Value = 38
$ws.Range("AJ42").Value = 11.5
$ws.Range("AM42").Value = 1000
$ws.Range("AO42").Value = 280
$ws.Range("F44").Value = 1.57
$ws.Range("I44").Value = 9.8
$ws.Range("K44").Value = 4.4
$ws.Range("T44").Value = 2.18
$ws.Range("U44").Value = 1.73
$ws.Range("AD44").Value = 38
$ws.Range("AE44").Value = 180
$ws.Range("AI44").Value = 170
$ws.Range("L45").Value = 1.51
$ws.Range("R45").Value = 1.24
$ws.Range("F46").Value = 1.89
$ws.Range("H46").Value = 4.9
$ws.Range("I46").Value = 5.4
$ws.Range("K46").Value = 3.6
$ws.Range("Q46").Value = 2.26
$ws.Range("R46").Value = 1.26
$ws.Range("S46").Value = 4.4
$ws.Range("V46").Value = 1.23
$ws.Range("X46").Value = 11.5
$ws.Range("Y46").Value = 15
$ws.Range("AI46").Value = 95
$ws.Range("F47").Value = 1.71
$ws.Range("G47").Value = 1.77
$ws.Range("H47").Value = 5.8
$ws.Range("I47").Value = 6.4
$ws.Range("K47").Value = 3.95
$ws.Range("T47").Value = 2.12
$ws.Range("AF47").Value = 9.6
$ws.Range("F48").Value = 2.26
$ws.Range("G48").Value = 2.4
$ws.Range("L48").Value = 1.59
$ws.Range("Q48").Value = 2.6
$ws.Range("R48").Value = 1.19
$ws.Range("W48").Value = 1.71
$ws.Range("AG48").Value = 12
$ws.Range("H49").Value = 2.52
$ws.Range("Q49").Value = 2.24
$ws.Range("F50").Value = 1.87
$ws.Range("H50").Value = 5.1
$ws.Range("J50").Value = 3.45
$ws.Range("T50").Value = 2.16
$ws.Range("AI50").Value = 120
$ws.Range("AO50").Value = 170
$ws.Range("AC51").Value = 10.5
$ws.Range("Q52").Value = 2.12
$ws.Range("M54").Value = 1.12
$ws.Range("N54").Value = 2.62
$ws.Range("O54").Value = 1.53
$ws.Range("P54").Value = 1.54
$ws.Range("Q54").Value = 2.54
$ws.Range("T54").Value = 2.06
$ws.Range("W54").Value = 1.3
$ws.Range("X54").Value = 10.5
$ws.Range("Y54").Value = 8.8
$ws.Range("Z54").Value = 16
$ws.Range("AA54").Value = 42
$ws.Range("AB54").Value = 13.5
$ws.Range("AC54").Value = 8.6
$ws.Range("AD54").Value = 14.5
$ws.Range("AE54").Value = 40
$ws.Range("AF54").Value = 34
$ws.Range("AG54").Value = 21
$ws.Range("AH54").Value = 27
$ws.Range("AI54").Value = 70
$ws.Range("AJ54").Value = 110
$ws.Range("AK54").Value = 80
$ws.Range("AL54").Value = 110
$ws.Range("AM54").Value = 210
$ws.Range("AN54").Value = 110
$ws.Range("AO54").Value = 40
$ws.Range("F55").Value = 1.45
$ws.Range("G55").Value = 1.52
$ws.Range("H55").Value = 9
$ws.Range("I55").Value = 11
$ws.Range("J55").Value = 4.1
$ws.Range("K55").Value = 4.6
$ws.Range("M55").Value = 1.08
$ws.Range("N55").Value = 3.2
$ws.Range("O55").Value = 1.38
$ws.Range("P55").Value = 1.73
$ws.Range("Q55").Value = 2.16
$ws.Range("R55").Value = 1.27
$ws.Range("S55").Value = 4
$ws.Range("T55").Value = 2.36
$ws.Range("U55").Value = 1.61
$ws.Range("V55").Value = 1.1
$ws.Range("W55").Value = 2.92
$ws.Range("X55").Value = 13.5
$ws.Range("Y55").Value = 25
$ws.Range("Z55").Value = 90
$ws.Range("AA55").Value = 480
$ws.Range("AB55").Value = 6.6
$ws.Range("AC55").Value = 10.5
$ws.Range("AD55").Value = 42
$ws.Range("AE55").Value = 250
$ws.Range("AF55").Value = 7.8
$ws.Range("AH55").Value = 34
$ws.Range("AI55").Value = 240
$ws.Range("AJ55").Value = 15
$ws.Range("AN55").Value = 11.5
$ws.Range("AO55").Value = 440
$ws.Range("H56").Value = 2
$ws.Range("J56").Value = 3.65
$ws.Range("N56").Value = 3.2
$ws.Range("P56").Value = 1.74
$ws.Range("S56").Value = 4.2
$ws.Range("U56").Value = 1.92
$ws.Range("Z56").Value = 13
$ws.Range("G57").Value = 2.88
$ws.Range("F58").Value = 3
$ws.Range("H58").Value = 2.68
$ws.Range("I58").Value = 2.9
$ws.Range("J58").Value = 3.05
$ws.Range("K58").Value = 3.25
$ws.Range("M58").Value = 1.12
$ws.Range("O58").Value = 1.56
$ws.Range("Q58").Value = 2.66
$ws.Range("V58").Value = 1.52
$ws.Range("W58").Value = 1.43
$ws.Range("AL58").Value = 80
